# Applies the "cryptos list" data refresh described in the commit
# "Updated cryptos list on Sun May 19 16:40:29 UTC 2024 with GitHub Actions".
# For each changed cell, set its new text value. Columns D sometimes contain
# plain numeric-looking strings (e.g. "575.90"); Excel would silently turn
# those into real numbers (and reformat them, e.g. "575.90" -> 575.9) unless
# the cell is forced to Text format first. Set-CellText does that, writes the
# value, then restores the default "Normal" style so no stray formatting is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '66.776.83'
$ws.Range('E2').Value = '  -0.20%  '

# Row 3
$ws.Range('D3').Value = '3.066.00'
$ws.Range('E3').Value = '  -1.55%  '

# Row 4
Set-CellText 'D4' '0.999'
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
Set-CellText 'D5' '575.90'
$ws.Range('E5').Value = '  -0.28%  '

# Row 6
Set-CellText 'D6' '168.12'
$ws.Range('E6').Value = '  -2.02%  '

# Row 7
Set-CellText 'D7' '0.999'
$ws.Range('E7').Value = '  -0.14%  '

# Row 8
$ws.Range('D8').Value = '3.063.63'
$ws.Range('E8').Value = '  -1.55%  '

# Row 9
$ws.Range('E9').Value = '  -2.01%  '

# Row 10
Set-CellText 'D10' '6.36'
$ws.Range('E10').Value = '  -1.63%  '

# Row 11
$ws.Range('E11').Value = '  -1.70%  '

# Row 12
Set-CellText 'D12' '0.468'
$ws.Range('E12').Value = '  -2.90%  '

# Row 13
Set-CellText 'D13' '0.0000239'
$ws.Range('E13').Value = '  -2.46%  '

# Row 14
Set-CellText 'D14' '35.62'
$ws.Range('E14').Value = '  -4.36%  '

# Row 15
$ws.Range('E15').Value = '  -1.86%  '

# Row 16
$ws.Range('D16').Value = '66.715.19'
$ws.Range('E16').Value = '  -0.24%  '

# Row 17
$ws.Range('D17').Value = '3.573.66'
$ws.Range('E17').Value = '  -1.65%  '

# Row 18
$ws.Range('E18').Value = '  -1.88%  '

# Row 19
Set-CellText 'D19' '16.76'
$ws.Range('E19').Value = '  +1.97%  '

# Row 20
$ws.Range('D20').Value = '3.061.53'
$ws.Range('E20').Value = '  -1.79%  '

# Row 21
Set-CellText 'D21' '489.86'
$ws.Range('E21').Value = '  +2.90%  '

# Row 22
$ws.Range('E22').Value = '  -3.45%  '

# Row 23
Set-CellText 'D23' '7.70'
$ws.Range('E23').Value = '  -3.21%  '

# Row 24
Set-CellText 'D24' '82.73'

# Row 25
Set-CellText 'D25' '12.68'
$ws.Range('E25').Value = '  -6.03%  '

# Row 26
Set-CellText 'D26' '2.20'
$ws.Range('E26').Value = '  -4.63%  '

# Row 27
Set-CellText 'D27' '10.17'
$ws.Range('E27').Value = '  +1.66%  '

# Row 28
$ws.Range('E28').Value = '  +0.03%  '

# Row 29
Set-CellText 'D29' '7.76'

# Row 30
Set-CellText 'D30' '2.26'
$ws.Range('E30').Value = '  -5.51%  '

# Row 31
$ws.Range('E31').Value = '  -2.38%  '

# Row 32
Set-CellText 'D32' '27.50'
$ws.Range('E32').Value = '  -3.62%  '

# Row 33
$ws.Range('E33').Value = '  -3.37%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0911'
$ws.Range('E34').Value = '  -2.58%  '

# Row 35
$ws.Range('E35').Value = '  +0.00%  '

# Row 36
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText 'D36' '0.948'
$ws.Range('E36').Value = '  -3.15%  '

# Row 37
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText 'D37' '5.60'
$ws.Range('E37').Value = '  -4.47%  '

# Row 38
Set-CellText 'D38' '46.57'
$ws.Range('E38').Value = '  -1.70%  '

# Row 39
$ws.Range('E39').Value = '  +0.24%  '

# Row 40
$ws.Range('E40').Value = '  -5.13%  '

# Row 41
$ws.Range('E41').Value = '  -3.51%  '

# Row 42
$ws.Range('E42').Value = '  -4.03%  '

# Row 43
$ws.Range('D43').Value = '2.754.12'
$ws.Range('E43').Value = '  -1.90%  '

# Row 44
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText 'D44' '0.0345'
$ws.Range('E44').Value = '  -3.20%  '

# Row 45
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText 'D45' '136.02'
$ws.Range('E45').Value = '  -0.01%  '

# Row 46
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-CellText 'D46' '365.87'
$ws.Range('E46').Value = '  -4.10%  '

# Row 47
Set-CellText 'D47' '2.47'
$ws.Range('E47').Value = '  -4.40%  '

# Row 48
$ws.Range('E48').Value = '  -0.04%  '

# Row 49
Set-CellText 'D49' '24.41'
$ws.Range('E49').Value = '  -1.06%  '

# Row 50
Set-CellText 'D50' '2.15'
$ws.Range('E50').Value = '  -2.18%  '

# Row 51
$ws.Range('E51').Value = '  -1.93%  '
